$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 46
$ws.Range("H46").Value = 1165.3334
$ws.Range("I46").Value = 998.5
$ws.Range("J46").Value = 1499
$ws.Range("K46").Value = 2995.5
$ws.Range("L46").Value = 4497
$ws.Range("M46").Value = -2876.5
$ws.Range("N46").Value = -4735

# Row 54
$ws.Range("H54").Value = 10000
$ws.Range("I54").Value = 5000
$ws.Range("K54").Value = 5000
$ws.Range("M54").Value = -4514

# Row 60
$ws.Range("H60").Value = 1165.3334
$ws.Range("I60").Value = 998.5
$ws.Range("J60").Value = 1499
$ws.Range("K60").Value = 2995.5
$ws.Range("L60").Value = 4497
$ws.Range("M60").Value = -2511.5
$ws.Range("N60").Value = -5465

# Row 74
$ws.Range("H74").Value = 6111.9585
$ws.Range("I74").Value = 3307.4443
$ws.Range("K74").Value = 3307.4443
$ws.Range("M74").Value = -2371.4443

# Row 77
$ws.Range("H77").Value = 6111.9585
$ws.Range("I77").Value = 3307.4443
$ws.Range("K77").Value = 16537.2215
$ws.Range("M77").Value = -11857.2215

# Row 82
$ws.Range("H82").Value = 3248.1667
$ws.Range("I82").Value = 3248.1667
$ws.Range("K82").Value = 9744.500100000001
$ws.Range("M82").Value = -9338.500100000001

# Row 85
$ws.Range("H85").Value = 3248.1667
$ws.Range("I85").Value = 3248.1667
$ws.Range("K85").Value = 9744.500100000001
$ws.Range("M85").Value = -8340.500100000001

# Row 137
$ws.Range("H137").Value = 2849.2632
$ws.Range("I137").Value = 1718.25
$ws.Range("K137").Value = 5154.75
$ws.Range("M137").Value = -2604.75

$ws = $wb.Worksheets.Item("ARM")
# Row 17
$ws.Range("H17").Value = 9999
$ws.Range("J17").Value = 9999
$ws.Range("L17").Value = 9999
$ws.Range("N17").Value = -10345

# Row 32
$ws.Range("H32").Value = 4273.3115
$ws.Range("I32").Value = 2785.2632
$ws.Range("K32").Value = 2785.2632
$ws.Range("M32").Value = -2498.2632

# Row 45
$ws.Range("H45").Value = 87626.25
$ws.Range("I45").Value = 144503.86
$ws.Range("J45").Value = 7997.6
$ws.Range("K45").Value = 144503.86
$ws.Range("L45").Value = 7997.6
$ws.Range("M45").Value = -144126.86
$ws.Range("N45").Value = -8751.6

# Row 74
$ws.Range("H74").Value = 539979.6
$ws.Range("I74").Value = 333666
$ws.Range("J74").Value = 849450
$ws.Range("K74").Value = 333666
$ws.Range("L74").Value = 849450
$ws.Range("M74").Value = -332792
$ws.Range("N74").Value = -851198

# Row 77
$ws.Range("H77").Value = 539979.6
$ws.Range("I77").Value = 333666
$ws.Range("J77").Value = 849450
$ws.Range("K77").Value = 1668330
$ws.Range("L77").Value = 4247250
$ws.Range("M77").Value = -1663962
$ws.Range("N77").Value = -4255986

# Row 122
$ws.Range("H122").Value = 3053.3667
$ws.Range("I122").Value = 2466.9524
$ws.Range("K122").Value = 7400.8572
$ws.Range("M122").Value = -4950.8572

$ws = $wb.Worksheets.Item("BSM")
# Row 13
$ws.Range("H13").Value = 53377.5
$ws.Range("J13").Value = 53377.5
$ws.Range("L13").Value = 53377.5
$ws.Range("N13").Value = -53713.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2617.549
$ws.Range("I31").Value = 2144.348
$ws.Range("J31").Value = 3006.25
$ws.Range("K31").Value = 2144.348
$ws.Range("L31").Value = 3006.25
$ws.Range("M31").Value = -1849.348
$ws.Range("N31").Value = -3596.25

# Row 34
$ws.Range("H34").Value = 2617.549
$ws.Range("I34").Value = 2144.348
$ws.Range("J34").Value = 3006.25
$ws.Range("K34").Value = 2144.348
$ws.Range("L34").Value = 3006.25
$ws.Range("M34").Value = -1942.348
$ws.Range("N34").Value = -3410.25

# Row 52
$ws.Range("H52").Value = 92963.336
$ws.Range("J52").Value = 92963.336
$ws.Range("L52").Value = 92963.336
$ws.Range("N52").Value = -93551.336

# Row 97
$ws.Range("H97").Value = 49715.57
$ws.Range("J97").Value = 49715.57
$ws.Range("L97").Value = 49715.57
$ws.Range("N97").Value = -51697.57

# Row 102
$ws.Range("H102").Value = 54993.668
$ws.Range("J102").Value = 54993.668
$ws.Range("L102").Value = 54993.668
$ws.Range("N102").Value = -59861.668

# Row 107
$ws.Range("H107").Value = 1606
$ws.Range("I107").Value = 1518.3889
$ws.Range("J107").Value = 1921.4
$ws.Range("K107").Value = 1518.3889
$ws.Range("L107").Value = 1921.4
$ws.Range("M107").Value = 401.6111000000001
$ws.Range("N107").Value = -5761.4

# Row 122
$ws.Range("H122").Value = 4438
$ws.Range("I122").Value = 4431.6665
$ws.Range("J122").Value = 4447.5
$ws.Range("K122").Value = 13294.9995
$ws.Range("L122").Value = 13342.5
$ws.Range("M122").Value = -10844.9995
$ws.Range("N122").Value = -18242.5

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 7448981
$ws.Range("I4").Value = 7486264
$ws.Range("K4").Value = 22458792
$ws.Range("M4").Value = -22458680

# Row 5
$ws.Range("H5").Value = 1619.5883
$ws.Range("I5").Value = 1153.7
$ws.Range("J5").Value = 2285.1428
$ws.Range("K5").Value = 3461.1
$ws.Range("L5").Value = 6855.428400000001
$ws.Range("M5").Value = -3349.1
$ws.Range("N5").Value = -7079.428400000001

# Row 63
$ws.Range("H63").Value = 13443.875
$ws.Range("I63").Value = 15167.4
$ws.Range("J63").Value = 10571.333
$ws.Range("K63").Value = 45502.2
$ws.Range("L63").Value = 31713.999
$ws.Range("M63").Value = -44753.2
$ws.Range("N63").Value = -33211.999

# Row 66
$ws.Range("H66").Value = 13443.875
$ws.Range("I66").Value = 15167.4
$ws.Range("J66").Value = 10571.333
$ws.Range("K66").Value = 136506.6
$ws.Range("L66").Value = 95141.997
$ws.Range("M66").Value = -132762.6
$ws.Range("N66").Value = -102629.997

# Row 131
$ws.Range("H131").Value = 13443052
$ws.Range("I131").Value = 7576416.5
$ws.Range("J131").Value = 16669702
$ws.Range("K131").Value = 22729249.5
$ws.Range("L131").Value = 50009106
$ws.Range("M131").Value = -22724209.5
$ws.Range("N131").Value = -50019186

# Row 135
$ws.Range("H135").Value = 1619.5883
$ws.Range("I135").Value = 1153.7
$ws.Range("J135").Value = 2285.1428
$ws.Range("K135").Value = 10383.3
$ws.Range("L135").Value = 20566.2852
$ws.Range("M135").Value = -7848.300000000001
$ws.Range("N135").Value = -25636.2852

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 699
$ws.Range("I22").Value = 764.4286
$ws.Range("J22").Value = 470
$ws.Range("K22").Value = 764.4286
$ws.Range("L22").Value = 470
$ws.Range("M22").Value = -469.4286
$ws.Range("N22").Value = -1060

# Row 27
$ws.Range("H27").Value = 699
$ws.Range("I27").Value = 764.4286
$ws.Range("J27").Value = 470
$ws.Range("K27").Value = 764.4286
$ws.Range("L27").Value = 470
$ws.Range("M27").Value = -657.4286
$ws.Range("N27").Value = -684

# Row 55
$ws.Range("H55").Value = 4723.727
$ws.Range("I55").Value = 4132.75
$ws.Range("K55").Value = 4132.75
$ws.Range("M55").Value = -3959.75

$ws = $wb.Worksheets.Item("WVR")
# Row 9
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

# Row 92
$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992

# Row 132
$ws.Range("H132").Value = 7316.44
$ws.Range("I132").Value = 8337.474
$ws.Range("K132").Value = 25012.422
$ws.Range("M132").Value = -22482.422

# Row 136
$ws.Range("H136").Value = 1756
$ws.Range("I136").Value = 1090.6511
$ws.Range("K136").Value = 3271.9533
$ws.Range("M136").Value = -721.9533000000001

Write-Host "Applied all Hyperion_Profits updates"
